# Update cryptos price/volume data and handle coin list shift (mCoin removed, USDD added)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.762.79'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.536.62'
$ws.Range("E3").Value = '  -1.95%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'205.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("D6").Value = "'0.485"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").Value = "'21.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.23%  '
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("D12").Value = '1.756.02'
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("D13").Value = '1.539.33'
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").Value = '26.764.47'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = "'60.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.04%  '
$ws.Range("D18").Value = "'212.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("E19").Value = '  -2.06%  '
$ws.Range("E20").Value = '  +0.52%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = "'4.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.31%  '
$ws.Range("D23").Value = "'9.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("E24").Value = '  -3.23%  '
$ws.Range("D25").Value = "'152.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("D26").Value = "'6.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.64%  '
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("E30").Value = '  -1.10%  '
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("E32").Value = '  +1.80%  '
$ws.Range("D33").Value = '1.365.76'
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("E34").Value = '  -0.38%  '
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("D36").Value = "'0.958"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.16%  '
$ws.Range("D37").Value = "'2.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("D39").Value = "'0.519"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.39%  '
$ws.Range("E40").Value = '  +7.62%  '
$ws.Range("E41").Value = '  -2.04%  '
$ws.Range("D42").Value = "'0.992"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").Value = "'62.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("D45").Value = "'1.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.49%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.670.59'
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = "'84.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = "'0.0507"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.16%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₇0980'
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = "'0.0942"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = "'1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.07%  '
